# Apply latest cryptos data update to sheet1 (Coin, Link, Price, Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''37.359.23'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('D3').Value = '''2.079.53'
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''250.74'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = '''0.662'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''56.75'
$ws.Range('E8').Value = '  +24.58%  '
$ws.Range('D9').Value = '''62.14'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Value = '''0.387'
$ws.Range('E10').Value = '  +4.87%  '
$ws.Range('D11').Value = '''0.0757'
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('E12').Value = '  +7.28%  '
$ws.Range('D13').Value = '''15.54'
$ws.Range('E13').Value = '  +5.91%  '
$ws.Range('E14').Value = '  -2.07%  '
$ws.Range('D15').Value = '''0.844'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = '''5.30'
$ws.Range('E16').Value = '  +3.84%  '
$ws.Range('D17').Value = '''2.081.66'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '''37.345.41'
$ws.Range('E18').Value = '  +1.57%  '
$ws.Range('D19').Value = '''73.23'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '''14.79'
$ws.Range('E20').Value = '  +13.83%  '
$ws.Range('D21').Value = '''0.0₃0851'
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('D22').Value = '''240.95'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '''5.28'
$ws.Range('E23').Value = '  +2.70%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = '''171.09'
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''21.07'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''9.21'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.12'
$ws.Range('E31').Value = '  +23.27%  '
$ws.Range('B32').Value = 'Gas'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D32').Value = '''23.17'
$ws.Range('E32').Value = '  +5.22%  '
$ws.Range('D33').Value = '''4.58'
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('D34').Value = '''0.0630'
$ws.Range('E34').Value = '  +5.05%  '
$ws.Range('E35').Value = '  +6.63%  '
$ws.Range('D36').Value = '''0.0912'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = '''1.86'
$ws.Range('E38').Value = '  -1.54%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '''2.29'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('E41').Value = '  +19.27%  '
$ws.Range('E42').Value = '  +3.77%  '
$ws.Range('D43').Value = '''17.77'
$ws.Range('E43').Value = '  +8.36%  '
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('D45').Value = '''100.12'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = '''4.30'
$ws.Range('E46').Value = '  +97.02%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value = '''2.80'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').Value = '''1.324.67'
$ws.Range('E48').Value = '  -3.41%  '
$ws.Range('D49').Value = '''2.42'
$ws.Range('E49').Value = '  +5.82%  '
$ws.Range('D50').Value = '''2.94'
$ws.Range('E50').Value = '  +3.34%  '
$ws.Range('D51').Value = '''7.03'
$ws.Range('E51').Value = '  +7.09%  '
